$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "46.725.04"
$ws.Range("E2").Value = "  +4.78%  "

# Row 3
$ws.Range("D3").Value = "2.331.91"
$ws.Range("E3").Value = "  +3.75%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.64%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.06%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.27%  "

# Row 7
$ws.Range("E7").Value = "  +0.87%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.58%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.88%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.59%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.46"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.19%  "

# Row 13
$ws.Range("E13").Value = "  -0.26%  "

# Row 14
$ws.Range("D14").Value = "2.689.68"
$ws.Range("E14").Value = "  +3.85%  "

# Row 15
$ws.Range("D15").Value = "2.337.66"
$ws.Range("E15").Value = "  +4.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.16"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.41%  "

# Row 17
$ws.Range("E17").Value = "  -0.17%  "

# Row 18
$ws.Range("D18").Value = "46.575.75"
$ws.Range("E18").Value = "  +4.91%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +10.31%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0946"
$ws.Range("E20").Value = "  +0.84%  "

# Row 21
$ws.Range("E21").Value = "  +0.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.88"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.37%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.33"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.35%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.52%  "

# Row 25
$ws.Range("E25").Value = "  -0.26%  "

# Row 26
$ws.Range("E26").Value = "  -0.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.09"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +12.93%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.72%  "

# Row 29
$ws.Range("E29").Value = "  +0.97%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.20"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.74"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.03%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.34"
$ws.Range("D32").ClearFormats()

# Row 33
$ws.Range("E33").Value = "  +2.71%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.62"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.48%  "

# Row 35
$ws.Range("E35").Value = "  -5.48%  "

# Row 36
$ws.Range("E36").Value = "  +0.55%  "

# Row 37
$ws.Range("E37").Value = "  -2.24%  "

# Row 38
$ws.Range("E38").Value = "  -3.38%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.42%  "

# Row 40
$ws.Range("E40").Value = "  +6.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.41"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.01"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.75%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.94"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.42%  "

# Row 45
$ws.Range("D45").Value = "1.788.93"
$ws.Range("E45").Value = "  -0.99%  "

# Row 46
$ws.Range("E46").Value = "  +6.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.55"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +8.51%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "80.76"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.27%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.66"
$ws.Range("D49").ClearFormats()

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.25"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.21%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.87"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.67%  "
